$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the release date on row 2 (44555.875 -> 44562.875)
$ws.Range("C2").Value = 44562.875

# Remove the extra data rows (old rows 3-6: FHIR/pos/cdh/java entries + the
# trailing blank row) - only one blank templated row remains afterwards.
$ws.Rows("3:6").Delete() | Out-Null

# Re-create the trailing blank row with the date number format preserved
# on column C (matches the workbook's original "template" row).
$ws.Range("C3").NumberFormat = "m/d/yy h:mm"

# Restore the selection to the new blank row.
$ws.Rows("3:3").Select() | Out-Null
